$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G3").Value = "2016-08-21 16:52:59"
$zhcn.Range("H3").Value = "2016-08-21 16:52:54"
$zhcn.Range("K3").Value = "2016-08-21 16:53:15"
$dede.Range("H3").Value = "2016-08-21 16:52:59"
$dede.Range("K3").Value = "2016-08-21 16:53:21"
